# Auto-generated script applying updated market price data to Sheets/Midgardsormr_Profits
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ
# columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR job sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 270.2
$ws.Cells.Item(4, 9).Value = 270.2
$ws.Cells.Item(4, 11).Value = 270.2
$ws.Cells.Item(4, 13).Value = -156.2
$ws.Cells.Item(11, 8).Value = 243.41667
$ws.Cells.Item(11, 9).Value = 243.41667
$ws.Cells.Item(11, 11).Value = 243.41667
$ws.Cells.Item(11, 13).Value = -103.41667
$ws.Cells.Item(17, 8).Value = 2778
$ws.Cells.Item(17, 10).Value = 3013.8
$ws.Cells.Item(17, 12).Value = 9041.400000000001
$ws.Cells.Item(17, 14).Value = -9377.400000000001
$ws.Cells.Item(40, 8).Value = 3200.1667
$ws.Cells.Item(40, 10).Value = 3825.25
$ws.Cells.Item(40, 12).Value = 3825.25
$ws.Cells.Item(40, 14).Value = -4175.25
$ws.Cells.Item(98, 8).Value = 6332
$ws.Cells.Item(98, 9).Value = 1999
$ws.Cells.Item(98, 11).Value = 1999
$ws.Cells.Item(98, 13).Value = -501
$ws.Cells.Item(113, 8).Value = 8771.611000000001
$ws.Cells.Item(113, 9).Value = 4558.7144
$ws.Cells.Item(113, 11).Value = 4558.7144
$ws.Cells.Item(113, 13).Value = -1304.7144
$ws.Cells.Item(116, 8).Value = 5820.923
$ws.Cells.Item(116, 10).Value = 5610.75
$ws.Cells.Item(116, 12).Value = 5610.75
$ws.Cells.Item(116, 14).Value = -12494.75
$ws.Cells.Item(122, 8).Value = 6332
$ws.Cells.Item(122, 9).Value = 1999
$ws.Cells.Item(122, 11).Value = 5997
$ws.Cells.Item(122, 13).Value = -3547
$ws.Cells.Item(131, 8).Value = 4502.143
$ws.Cells.Item(131, 9).Value = 4091.5
$ws.Cells.Item(131, 10).Value = 6966
$ws.Cells.Item(131, 11).Value = 12274.5
$ws.Cells.Item(131, 12).Value = 20898
$ws.Cells.Item(131, 13).Value = -7234.5
$ws.Cells.Item(131, 14).Value = -30978
$ws.Cells.Item(132, 8).Value = 25124.678
$ws.Cells.Item(132, 9).Value = 27281.285
$ws.Cells.Item(132, 11).Value = 81843.855
$ws.Cells.Item(132, 13).Value = -79313.855
$ws.Cells.Item(135, 8).Value = 17522.361
$ws.Cells.Item(135, 9).Value = 1398.8
$ws.Cells.Item(135, 10).Value = 54166.816
$ws.Cells.Item(135, 11).Value = 12589.2
$ws.Cells.Item(135, 12).Value = 487501.344
$ws.Cells.Item(135, 13).Value = -10054.2
$ws.Cells.Item(135, 14).Value = -492571.344
$ws.Cells.Item(141, 8).Value = 1126.4286
$ws.Cells.Item(141, 9).Value = 1126.4286
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 3379.2858
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = 1800.7142
$ws.Cells.Item(141, 14).ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1385.9117
$ws.Cells.Item(2, 9).Value = 950.96
$ws.Cells.Item(2, 11).Value = 950.96
$ws.Cells.Item(2, 13).Value = -837.96
$ws.Cells.Item(32, 8).Value = 37185.44
$ws.Cells.Item(32, 9).Value = 40838.137
$ws.Cells.Item(32, 11).Value = 40838.137
$ws.Cells.Item(32, 13).Value = -40551.137
$ws.Cells.Item(36, 8).Value = 10724.8
$ws.Cells.Item(36, 9).Value = 4575.3335
$ws.Cells.Item(36, 11).Value = 4575.3335
$ws.Cells.Item(36, 13).Value = -4229.3335
$ws.Cells.Item(37, 8).Value = 25833.166
$ws.Cells.Item(37, 10).Value = 47499.5
$ws.Cells.Item(37, 12).Value = 47499.5
$ws.Cells.Item(37, 14).Value = -48045.5
$ws.Cells.Item(44, 8).Value = 57999
$ws.Cells.Item(44, 10).Value = 57999
$ws.Cells.Item(44, 12).Value = 57999
$ws.Cells.Item(44, 14).Value = -58975
$ws.Cells.Item(55, 8).Value = 47999
$ws.Cells.Item(55, 10).Value = 47999
$ws.Cells.Item(55, 12).Value = 47999
$ws.Cells.Item(55, 14).Value = -48629
$ws.Cells.Item(61, 8).Value = 2389.4075
$ws.Cells.Item(61, 9).Value = 1438.125
$ws.Cells.Item(61, 10).Value = 9999.666999999999
$ws.Cells.Item(61, 11).Value = 1438.125
$ws.Cells.Item(61, 12).Value = 9999.666999999999
$ws.Cells.Item(61, 13).Value = -1226.125
$ws.Cells.Item(61, 14).Value = -10423.667
$ws.Cells.Item(97, 8).Value = 1947.037
$ws.Cells.Item(97, 9).Value = 1393.4445
$ws.Cells.Item(97, 11).Value = 1393.4445
$ws.Cells.Item(97, 13).Value = -897.4445000000001
$ws.Cells.Item(116, 8).Value = 1385.9117
$ws.Cells.Item(116, 9).Value = 950.96
$ws.Cells.Item(116, 11).Value = 950.96
$ws.Cells.Item(116, 13).Value = 1343.04
$ws.Cells.Item(122, 8).Value = 2441.8096
$ws.Cells.Item(122, 9).Value = 2363.9
$ws.Cells.Item(122, 11).Value = 7091.700000000001
$ws.Cells.Item(122, 13).Value = -4641.700000000001
$ws.Cells.Item(136, 8).Value = 2389.4075
$ws.Cells.Item(136, 9).Value = 1438.125
$ws.Cells.Item(136, 10).Value = 9999.666999999999
$ws.Cells.Item(136, 11).Value = 4314.375
$ws.Cells.Item(136, 12).Value = 29999.001
$ws.Cells.Item(136, 13).Value = -1764.375
$ws.Cells.Item(136, 14).Value = -35099.001
$ws.Cells.Item(139, 8).Value = 149997.5
$ws.Cells.Item(139, 10).Value = 149997.5
$ws.Cells.Item(139, 12).Value = 149997.5
$ws.Cells.Item(139, 14).Value = -160277.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1385.9117
$ws.Cells.Item(3, 9).Value = 950.96
$ws.Cells.Item(3, 11).Value = 950.96
$ws.Cells.Item(3, 13).Value = -836.96
$ws.Cells.Item(138, 8).Value = 302926
$ws.Cells.Item(138, 10).Value = 302926
$ws.Cells.Item(138, 12).Value = 302926
$ws.Cells.Item(138, 14).Value = -313206

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 22226.6
$ws.Cells.Item(41, 10).Value = 36271.445
$ws.Cells.Item(41, 12).Value = 36271.445
$ws.Cells.Item(41, 14).Value = -37127.445
$ws.Cells.Item(50, 8).Value = 34498.223
$ws.Cells.Item(50, 10).Value = 34498.223
$ws.Cells.Item(50, 12).Value = 34498.223
$ws.Cells.Item(50, 14).Value = -35748.223
$ws.Cells.Item(51, 8).Value = 36979.2
$ws.Cells.Item(51, 10).Value = 36979.2
$ws.Cells.Item(51, 12).Value = 36979.2
$ws.Cells.Item(51, 14).Value = -38451.2
$ws.Cells.Item(59, 8).Value = 26574.8
$ws.Cells.Item(59, 10).Value = 26969
$ws.Cells.Item(59, 12).Value = 26969
$ws.Cells.Item(59, 14).Value = -29259
$ws.Cells.Item(60, 8).Value = 33288.168
$ws.Cells.Item(60, 10).Value = 33046
$ws.Cells.Item(60, 12).Value = 33046
$ws.Cells.Item(60, 14).Value = -34068
$ws.Cells.Item(61, 8).Value = 36979.2
$ws.Cells.Item(61, 10).Value = 36979.2
$ws.Cells.Item(61, 12).Value = 36979.2
$ws.Cells.Item(61, 14).Value = -37675.2
$ws.Cells.Item(94, 8).Value = 2500
$ws.Cells.Item(94, 10).Value = 2500
$ws.Cells.Item(94, 12).Value = 2500
$ws.Cells.Item(94, 14).Value = -3402
$ws.Cells.Item(99, 8).Value = 11367.538
$ws.Cells.Item(99, 9).Value = 10358
$ws.Cells.Item(99, 11).Value = 10358
$ws.Cells.Item(99, 13).Value = -8860
$ws.Cells.Item(126, 8).Value = 11367.538
$ws.Cells.Item(126, 9).Value = 10358
$ws.Cells.Item(126, 11).Value = 31074
$ws.Cells.Item(126, 13).Value = -28604
$ws.Cells.Item(132, 8).Value = 101238.664
$ws.Cells.Item(132, 9).Value = 151133
$ws.Cells.Item(132, 10).Value = 1450
$ws.Cells.Item(132, 11).Value = 453399
$ws.Cells.Item(132, 12).Value = 4350
$ws.Cells.Item(132, 13).Value = -450869
$ws.Cells.Item(132, 14).Value = -9410
$ws.Cells.Item(134, 8).Value = 2114.6487
$ws.Cells.Item(134, 9).Value = 1797.0625
$ws.Cells.Item(134, 11).Value = 5391.1875
$ws.Cells.Item(134, 13).Value = -2856.1875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 270
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(33, 8).Value = 89
$ws.Cells.Item(33, 9).Value = 90.8
$ws.Cells.Item(33, 10).Value = 80
$ws.Cells.Item(33, 11).Value = 544.8
$ws.Cells.Item(33, 12).Value = 480
$ws.Cells.Item(33, 13).Value = -261.8
$ws.Cells.Item(33, 14).Value = -1046
$ws.Cells.Item(44, 8).Value = 2749.5
$ws.Cells.Item(44, 10).Value = 2499
$ws.Cells.Item(44, 12).Value = 7497
$ws.Cells.Item(44, 14).Value = -8293

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 6604.9
$ws.Cells.Item(43, 9).Value = 873.2857
$ws.Cells.Item(43, 10).Value = 19978.666
$ws.Cells.Item(43, 11).Value = 873.2857
$ws.Cells.Item(43, 12).Value = 19978.666
$ws.Cells.Item(43, 13).Value = -722.2857
$ws.Cells.Item(43, 14).Value = -20280.666
$ws.Cells.Item(122, 8).Value = 5216.5835
$ws.Cells.Item(122, 9).Value = 5233.222
$ws.Cells.Item(122, 10).Value = 5166.6665
$ws.Cells.Item(122, 11).Value = 15699.666
$ws.Cells.Item(122, 12).Value = 15499.9995
$ws.Cells.Item(122, 13).Value = -13249.666
$ws.Cells.Item(122, 14).Value = -20399.9995
$ws.Cells.Item(123, 8).Value = 74999.5
$ws.Cells.Item(123, 10).Value = 74999.5
$ws.Cells.Item(123, 12).Value = 74999.5
$ws.Cells.Item(123, 14).Value = -79899.5
$ws.Cells.Item(132, 8).Value = 3659.2856
$ws.Cells.Item(132, 9).Value = 2644.7812
$ws.Cells.Item(132, 11).Value = 7934.3436
$ws.Cells.Item(132, 13).Value = -5404.3436

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 68692560
$ws.Cells.Item(12, 9).Value = 68692560
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 68692560
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -68692390
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 3610.3635
$ws.Cells.Item(40, 9).Value = 3557.647
$ws.Cells.Item(40, 11).Value = 3557.647
$ws.Cells.Item(40, 13).Value = -3421.647
$ws.Cells.Item(136, 8).Value = 6164.077
$ws.Cells.Item(136, 9).Value = 6241.864
$ws.Cells.Item(136, 11).Value = 18725.592
$ws.Cells.Item(136, 13).Value = -16175.592

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 826.7778
$ws.Cells.Item(107, 9).Value = 826.7778
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2480.3334
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -560.3334
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 39256.94
$ws.Cells.Item(122, 9).Value = 40969.484
$ws.Cells.Item(122, 11).Value = 122908.452
$ws.Cells.Item(122, 13).Value = -120458.452
$ws.Cells.Item(132, 8).Value = 2060.0222
$ws.Cells.Item(132, 9).Value = 2132.7856
$ws.Cells.Item(132, 11).Value = 6398.3568
$ws.Cells.Item(132, 13).Value = -3868.3568
$ws.Cells.Item(136, 8).Value = 21374.205
$ws.Cells.Item(136, 9).Value = 22397.656
$ws.Cells.Item(136, 11).Value = 67192.96799999999
$ws.Cells.Item(136, 13).Value = -64642.96799999999
